$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill the duplicate_image_filename column (E) with "NA" for data rows 2-21
$ws.Range("E2:E21").Value = "NA"
